$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the "Discharged" / "Intubations" column headers (D1 <-> E1) ---
$ws.Range("D1").Value = "Intubations"
$ws.Range("E1").Value = "Discharged"

# --- Swap the D/E values for every existing data row (2-21) ---
# Row 2 only had a D value (E2 was empty) -> move it to E2 and clear D2.
$ws.Range("E2").Value = 0
$ws.Range("D2").ClearContents()

$swapRows = @{
    3  = @(0, 43)
    4  = @(123, 27)
    5  = @(78, 93)
    6  = @(102, 69)
    7  = @(150, 116)
    8  = @(145, 124)
    9  = @(150, 121)
    10 = @(278, 192)
    11 = @(450, 222)
    12 = @(528, 290)
    13 = @(681, 200)
    14 = @(846, 165)
    15 = @(632, 303)
    16 = @(771, 295)
    17 = @(1167, 291)
    18 = @(1292, 313)
    19 = @(1452, 260)
    20 = @(1592, 351)
    21 = @(1709, 316)
}

foreach ($r in $swapRows.Keys) {
    $oldD = $swapRows[$r][0]
    $oldE = $swapRows[$r][1]
    $ws.Cells.Item($r, 4).Value = $oldE
    $ws.Cells.Item($r, 5).Value = $oldD
}

# --- Append the two new rows of data (5 April and 6 April 2020) ---
$ws.Range("A22").Value = 43926
$ws.Range("B22").Value = 358
$ws.Range("C22").Value = 128
$ws.Range("D22").Value = 132
$ws.Range("E22").Value = 1179
$ws.Range("A22").Style = $ws.Range("A21").Style

$ws.Range("A23").Value = 43927
$ws.Range("B23").Value = 656
$ws.Range("C23").Value = 89
$ws.Range("D23").Value = 69
$ws.Range("E23").Value = 1224
$ws.Range("A23").Style = $ws.Range("A21").Style

# --- Move the active selection to mirror the author's final cursor position ---
$ws.Range("E24").Select()
